$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sunscreens")

# Column A header: "Version" -> "ID"
$ws.Range("A1").Value = "ID"

# Renumber the ID column for existing rows and fill in the two new rows
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Update the Image column to point at the new Images/ folder naming scheme
$ws.Range("X2").Value = "Images/0001.png"
$ws.Range("X3").Value = "Images/0002.png"
$ws.Range("X4").Value = "Images/0003.png"
$ws.Range("X5").Value = "Images/0004.png"
$ws.Range("X6").Value = "Images/0005.png"

# Match the saved selection state
[void]$ws.Range("X7").Select()
